$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (distance, MEAN, STD, MIN, MAX, COUNT, Month)
$newRows = @(
    @(1100, 36.45192337036133, 3.148383855819702, 19.09218978881836, 48.14522552490234, 18476, "05"),
    @(1200, 36.44610214233398, 3.229307174682617, 18.54372787475586, 47.85818099975586, 18568, "05"),
    @(1300, 36.37449645996094, 3.262782573699951, 18.85298538208008, 46.3990364074707,  18438, "05"),
    @(1400, 36.27005004882812, 3.278586387634277, 18.97258758544922, 46.96799850463867, 18366, "05"),
    @(1500, 36.18234252929688, 3.3784499168396,   18.52151679992676, 50.05373001098633, 18392, "05")
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    # Column G holds a text "Month" code like "05" - force text format so
    # Excel doesn't auto-coerce the numeric-looking string into a number.
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row[6]
}
